$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths for the three date/time columns (B, C, H)
$ws.Columns(2).ColumnWidth = 20.7109375
$ws.Columns(3).ColumnWidth = 20.7109375
$ws.Columns(8).ColumnWidth = 20.7109375

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 44523.82420138889
$ws.Range("C2").Value = 44523.82730324075
$ws.Range("D2").Value = "IP Address"
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 267
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = 44523.82730324075
$ws.Range("I2").Value = "1bimil"
$ws.Range("J2").Value = "ebola %>% `n  pivot_longer(Cases_Guinea:last_col()) %>% `n  separate(name, into = c(`"case_death`", `"country`"), sep = `"_`") %>% `n  drop_na()`n"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 44523.10084490741
$ws.Range("C3").Value = 44523.10333333333
$ws.Range("D3").Value = "Spam"
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 215
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 44523.87929398148
$ws.Range("I3").Value = "2nesch"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 44523.82427083334
$ws.Range("C4").Value = 44523.82444444444
$ws.Range("D4").Value = "IP Address"
$ws.Range("E4").Value = 50
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 44523.87930555556
$ws.Range("I4").Value = "1lacat"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 44518.816875
$ws.Range("C5").Value = 44518.81747685185
$ws.Range("D5").Value = "IP Address"
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = 52
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 44523.87930555556
$ws.Range("I5").Value = "6chhom"

# Date/time formatting for columns B, C, H on the data rows
$ws.Range("B2:B5").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("C2:C5").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("H2:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"

Write-Output "done"
